$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(218).Insert()

$ws.Cells.Item(218, 1).Value = 9
$ws.Cells.Item(218, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(218, 3).Value = "Metropolitana"
$ws.Cells.Item(218, 4).Value = 44784
$ws.Cells.Item(218, 5).Value = 13
$ws.Cells.Item(218, 6).Value = 300000001
$ws.Cells.Item(218, 7).Value = "Rabanito"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 7000
$ws.Cells.Item(218, 11).Value = 2500
$ws.Cells.Item(218, 12).Value = 3000
$ws.Cells.Item(218, 13).Value = 2750
$ws.Cells.Item(218, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(218, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(218, 16).Value = 28
$ws.Cells.Item(218, 17).Value = 100
$ws.Cells.Item(218, 18).Value = "Hortaliza"
